$d = $word.ActiveDocument

# The hidden "_GoBack" bookmark currently sits at the end of the
# "Kevin can now get the e02 file from the service server." paragraph.
# It needs to move to the end of the new "...worked on the UI part"
# paragraph being added below, so delete it here first; it is re-created
# (embedded directly in the inserted XML) at its new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Insert the new daily-report entries (2nd June & 3rd June sections)
# right after the last paragraph, before the section break.
$insertXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Overskrift2"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>The 2</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/><w:lang w:val="en-US"/></w:rPr><w:t>nd</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> June</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">We scaped a lot, can not remember what we made. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Kristian and Asger worked on the server and office connection so that Asger could run the server and office on his PC and upload the e02 file to the server on his own Laptop. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Kevin worked on the Design background on the Orderoverview and worked on the UI part </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">They had some problems with the wrapper. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Overskrift2"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>The 3</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/><w:lang w:val="en-US"/></w:rPr><w:t>rd</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> June</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Asger made the “Order Overview” GUI for the office, had some problems with List view and well did some hard coding instead</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Kristian redesigned the SQL, rewrote the Datahandler and made some server fixes. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Kevin made some Android, and just need to test it with some data from the Database, but is not waiting for Kristian</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The Order parser is now ready and able to read their newest e02 file. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r = $d.Content
$r.Collapse(0)
$r.InsertXML($insertXml)

Write-Host "Paragraphs after insert:" $d.Paragraphs.Count
